# This script updates the NC initial-conditions table on Sheet1 to reflect
# the new data-processing run: row labels are renumbered (NC-MFR-ABS-1_4 ..
# NC-MFR-ABS-10_4, then NC-MFR-ABS-4 / NC-MFR_HE-ABS-4), and the numeric
# columns (B..J) are refreshed with the newly computed values. The selected
# range on the sheet is also updated to reflect the new active data block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update data rows 2-13 on Sheet1 (NC initial conditions table)
$ws.Cells.Item(2, 1).Value = "NC-MFR-ABS-1_4"
$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 4
$ws.Cells.Item(2, 7).Value = 143.6125329983828
$ws.Cells.Item(2, 8).Value = 123.61253299838279
$ws.Cells.Item(2, 9).Value = 2
$ws.Cells.Item(2, 10).Value = 20

$ws.Cells.Item(3, 1).Value = "NC-MFR-ABS-2_4"
$ws.Cells.Item(3, 2).Value = 0.20890589895780506
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = 0.05
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 4
$ws.Cells.Item(3, 7).Value = 141.77262336984768
$ws.Cells.Item(3, 8).Value = 121.77262336984768
$ws.Cells.Item(3, 9).Value = 2
$ws.Cells.Item(3, 10).Value = 20

$ws.Cells.Item(4, 1).Value = "NC-MFR-ABS-3_4"
$ws.Cells.Item(4, 2).Value = 0.3961441459088963
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = 0.1
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 4
$ws.Cells.Item(4, 7).Value = 139.85299129047138
$ws.Cells.Item(4, 8).Value = 119.85299129047138
$ws.Cells.Item(4, 9).Value = 2
$ws.Cells.Item(4, 10).Value = 20

$ws.Cells.Item(5, 1).Value = "NC-MFR-ABS-4_4"
$ws.Cells.Item(5, 2).Value = 0.5852513415477354
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(5, 4).Value = 0.15
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 4
$ws.Cells.Item(5, 7).Value = 137.84534219941258
$ws.Cells.Item(5, 8).Value = 117.84534219941258
$ws.Cells.Item(5, 9).Value = 2
$ws.Cells.Item(5, 10).Value = 20

$ws.Cells.Item(6, 1).Value = "NC-MFR-ABS-5_4"
$ws.Cells.Item(6, 2).Value = 0.7763949458931612
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(6, 4).Value = 0.2
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 4
$ws.Cells.Item(6, 7).Value = 135.73996041273483
$ws.Cells.Item(6, 8).Value = 115.73996041273483
$ws.Cells.Item(6, 9).Value = 2
$ws.Cells.Item(6, 10).Value = 20

$ws.Cells.Item(7, 1).Value = "NC-MFR-ABS-6_4"
$ws.Cells.Item(7, 2).Value = 0.9697714644966886
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 0.25
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 4
$ws.Cells.Item(7, 7).Value = 133.52535794654545
$ws.Cells.Item(7, 8).Value = 113.52535794654545
$ws.Cells.Item(7, 9).Value = 2
$ws.Cells.Item(7, 10).Value = 20

$ws.Cells.Item(8, 1).Value = "NC-MFR-ABS-7_4"
$ws.Cells.Item(8, 2).Value = 1.1656126975459682
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(8, 4).Value = 0.3
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 4
$ws.Cells.Item(8, 7).Value = 131.18780629572507
$ws.Cells.Item(8, 8).Value = 111.18780629572507
$ws.Cells.Item(8, 9).Value = 2
$ws.Cells.Item(8, 10).Value = 20

$ws.Cells.Item(9, 1).Value = "NC-MFR-ABS-8_4"
$ws.Cells.Item(9, 2).Value = 1.565862322646411
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = 0.4
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 4
$ws.Cells.Item(9, 7).Value = 126.07367775380408
$ws.Cells.Item(9, 8).Value = 106.07367775380408
$ws.Cells.Item(9, 9).Value = 2
$ws.Cells.Item(9, 10).Value = 20

$ws.Cells.Item(10, 1).Value = "NC-MFR-ABS-9_4"
$ws.Cells.Item(10, 2).Value = 1.9802133138045936
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(10, 4).Value = 0.5
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 4
$ws.Cells.Item(10, 7).Value = 120.21154593648862
$ws.Cells.Item(10, 8).Value = 100.21154593648862
$ws.Cells.Item(10, 9).Value = 2
$ws.Cells.Item(10, 10).Value = 20

$ws.Cells.Item(11, 1).Value = "NC-MFR-ABS-10_4"
$ws.Cells.Item(11, 2).Value = 2.639892051933314
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = 0.65
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 4
$ws.Cells.Item(11, 7).Value = 109.29210588198293
$ws.Cells.Item(11, 8).Value = 89.29210588198293
$ws.Cells.Item(11, 9).Value = 2
$ws.Cells.Item(11, 10).Value = 20

$ws.Cells.Item(12, 1).Value = "NC-MFR-ABS-4"
$ws.Cells.Item(12, 2).Value = 1.4245782176988584
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(12, 4).Value = 0.3
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 5
$ws.Cells.Item(12, 7).Value = 138.8607391541683
$ws.Cells.Item(12, 8).Value = 118.8607391541683
$ws.Cells.Item(12, 9).Value = 2
$ws.Cells.Item(12, 10).Value = 20

$ws.Cells.Item(13, 1).Value = "NC-MFR_HE-ABS-4"
$ws.Cells.Item(13, 2).Value = 1.4245782176988584
$ws.Cells.Item(13, 3).Value = 1.4245782176988584
$ws.Cells.Item(13, 4).Value = 0.3
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 5
$ws.Cells.Item(13, 7).Value = 138.8607391541683
$ws.Cells.Item(13, 8).Value = 118.8607391541683
$ws.Cells.Item(13, 9).Value = 2
$ws.Cells.Item(13, 10).Value = 20

# Update the visible selection to cover the new active data block (A1:J11)
$ws.Range("A1:J11").Select()
